{"js": "// Add the new \"TheoremStyleUpright\" paragraph style, based on the existing\n// \"TheoremStyle\" style, matching the author's change to word/styles.xml:\n//\n//   <w:style w:type=\"paragraph\" w:customStyle=\"1\" w:styleId=\"TheoremStyleUpright\">\n//     <w:name w:val=\"TheoremStyleUpright\"/>\n//     <w:basedOn w:val=\"TheoremStyle\"/>\n//     <w:qFormat/>\n//     <w:rsid w:val=\"009010A4\"/>\n//   </w:style>\n//\n// NOTE: `addStyle` returns a proxy anchored to the request that will create\n// the style; it is not safe to keep setting properties on that same proxy\n// before the style is materialized server-side. We sync once so the style\n// exists, then re-fetch it by name before setting its remaining properties.\nconst newStyle = context.document.addStyle(\"TheoremStyleUpright\", Word.StyleType.paragraph);\nawait context.sync();\n\nconst styles = context.document.getStyles();\nconst theoremStyleUpright = styles.getByNameOrNullObject(\"TheoremStyleUpright\");\nawait context.sync();\n\nif (!theoremStyleUpright.isNullObject) {\n  theoremStyleUpright.baseStyle = \"TheoremStyle\";\n  theoremStyleUpright.quickStyle = true;\n  await context.sync();\n}\n", "ps1": "# Add the new \"TheoremStyleUpright\" paragraph style, based on the existing\n# \"TheoremStyle\" style, matching the author's change to word/styles.xml.\n$d = $word.ActiveDocument\n\n$newStyle = $d.Styles.Add(\"TheoremStyleUpright\", 1)\n$newStyle.BaseStyle = $d.Styles.Item(\"TheoremStyle\")\n$newStyle.QuickStyle = $true\n"}
